# Rewrite the "Questions" sheet content:
#  - A1 previously held a bold/bordered numeric placeholder (0); it now
#    receives the (reformatted) question-bank text that used to live in A2.
#  - A2 is cleared out entirely (its old shared-string reference is dropped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @"
questions = [
    {
        "title": "As a solutions architect for a company operating multiple applications across various virtual private clouds (VPCs) and AWS accounts, your task is to design a centralized network architecture. This architecture should interconnect these VPCs and ensure that traffic to AWS does not traverse the public Internet.Which combination of AWS services should you use?",
        "ques_type": 2,
        "options": [
            "AWS Transit Gateway and VPC endpoints",
            "Amazon VPC Peering and AWS Direct Connect",
            "AWS Site-to-Site VPN and VPC endpoints",
            "AWS Transit Gateway and Amazon Route 53 Resolver"
        ],
        "score": "AWS Transit Gateway and VPC endpoints"
    },
    {
        "title": "As a cloud security specialist for an e-commerce platform, you are tasked with ensuring the platform's stability during peak sales periods. The platform experiences significant web traffic, including malicious attempts to disrupt service. You need a solution to protect against web threats and handle high traffic volumes without compromising performance.Which combination of AWS services should you use?",
        "ques_type": 2,
        "options": [
            "AWS WAF and Amazon CloudFront",
            "AWS Shield Advanced and Amazon Inspector",
            "AWS WAF and AWS VPN",
            "Amazon GuardDuty and AWS Shield Standard"
        ],
        "score": "AWS WAF and Amazon CloudFront"
    },
    {
        "title": "As a network security specialist, you have detected irregularities in the network traffic associated with one of your application's virtual private clouds (VPCs). To diagnose any potential security or configuration issues, you need in-depth visibility into the IP traffic to and from the network interfaces. Which AWS feature should you enable?",
        "ques_type": 2,
        "options": [
            "VPC Flow Logs",
            "Amazon CloudWatch Alarms",
            "AWS CloudTrail Logging",
            "AWS X-Ray Tracing"
        ],
        "score": "VPC Flow Logs"
    },
    {
        "title": "As a cloud infrastructure specialist working for a global media company, you are dealing with complaints about buffering issues from international viewers. Your goal is to distribute traffic over the AWS global network to ensure consistent application performance for all users.Which AWS service should you use?",
        "ques_type": 2,
        "options": [
            "Amazon CloudFront",
            "AWS Direct Connect",
            "AWS Global Accelerator",
            "AWS Transit Gateway"
        ],
        "score": "Amazon CloudFront"
    }
]
"@

# Drop the old A2 cell completely (no more second row).
$ws.Range("A2").ClearContents()

# Reset A1's formatting first (drops the bold font / border / centered
# alignment it had as the placeholder "0" cell) and only then write the
# long text, so the autofit pass below sizes a plain-format cell.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $newText

# Re-run the row's autofit so the embedded newlines in the new text don't
# leave a stale custom row height behind.
$ws.Range("A1").EntireRow.AutoFit()
